$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Btc"
$ws.Cells.Item(2, 3).Value = "Egfr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.191714
$ws.Cells.Item(2, 8).Value = 0.575142
$ws.Cells.Item(2, 9).Value = 0.09369188973541917
$ws.Cells.Item(2, 10).Value = 0.09369188973541917
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.701354
$ws.Cells.Item(2, 14).Value = 8.104061999999999
$ws.Cells.Item(2, 15).Value = 0.02221077311549548
$ws.Cells.Item(2, 16).Value = 0.02221077311549548
$ws.Cells.Item(2, 17).Value = 0.517887380756
$ws.Cells.Item(2, 18).Value = 4.660986426804
$ws.Cells.Item(2, 19).Value = 0.002080969305675415
$ws.Cells.Item(2, 20).Value = 0.002080969305675415

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Btc"
$ws.Cells.Item(3, 3).Value = "Egfr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.191714
$ws.Cells.Item(3, 8).Value = 0.575142
$ws.Cells.Item(3, 9).Value = 0.09369188973541917
$ws.Cells.Item(3, 10).Value = 0.09369188973541917
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 88.14978533333333
$ws.Cells.Item(3, 14).Value = 264.449356
$ws.Cells.Item(3, 15).Value = 0.7247753838328104
$ws.Cells.Item(3, 16).Value = 0.7247753838328105
$ws.Cells.Item(3, 17).Value = 16.89954794539467
$ws.Cells.Item(3, 18).Value = 152.095931508552
$ws.Cells.Item(3, 19).Value = 0.06790557534500978
$ws.Cells.Item(3, 20).Value = 0.06790557534500978

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Btc"
$ws.Cells.Item(4, 3).Value = "Egfr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.191714
$ws.Cells.Item(4, 8).Value = 0.575142
$ws.Cells.Item(4, 9).Value = 0.09369188973541917
$ws.Cells.Item(4, 10).Value = 0.09369188973541917
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.24063
$ws.Cells.Item(4, 14).Value = 0.72189
$ws.Cells.Item(4, 15).Value = 0.001978481285600361
$ws.Cells.Item(4, 16).Value = 0.001978481285600361
$ws.Cells.Item(4, 17).Value = 0.04613213982000001
$ws.Cells.Item(4, 18).Value = 0.41518925838
$ws.Cells.Item(4, 19).Value = 0.0001853676504540594
$ws.Cells.Item(4, 20).Value = 0.0001853676504540594

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Btc"
$ws.Cells.Item(5, 3).Value = "Egfr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.191714
$ws.Cells.Item(5, 8).Value = 0.575142
$ws.Cells.Item(5, 9).Value = 0.09369188973541917
$ws.Cells.Item(5, 10).Value = 0.09369188973541917
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 30.53182233333333
$ws.Cells.Item(5, 14).Value = 91.595467
$ws.Cells.Item(5, 15).Value = 0.2510353617660938
$ws.Cells.Item(5, 16).Value = 0.2510353617660938
$ws.Cells.Item(5, 17).Value = 5.853377786812668
$ws.Cells.Item(5, 18).Value = 52.680400081314
$ws.Cells.Item(5, 19).Value = 0.02351997743427992
$ws.Cells.Item(5, 20).Value = 0.02351997743427992

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Btc"
$ws.Cells.Item(6, 3).Value = "Egfr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.854503666666667
$ws.Cells.Item(6, 8).Value = 5.563511
$ws.Cells.Item(6, 9).Value = 0.9063081102645809
$ws.Cells.Item(6, 10).Value = 0.9063081102645809
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.701354
$ws.Cells.Item(6, 14).Value = 8.104061999999999
$ws.Cells.Item(6, 15).Value = 0.02221077311549548
$ws.Cells.Item(6, 16).Value = 0.02221077311549548
$ws.Cells.Item(6, 17).Value = 5.009670897964667
$ws.Cells.Item(6, 18).Value = 45.087038081682
$ws.Cells.Item(6, 19).Value = 0.02012980380982007
$ws.Cells.Item(6, 20).Value = 0.02012980380982007

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Btc"
$ws.Cells.Item(7, 3).Value = "Egfr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.854503666666667
$ws.Cells.Item(7, 8).Value = 5.563511
$ws.Cells.Item(7, 9).Value = 0.9063081102645809
$ws.Cells.Item(7, 10).Value = 0.9063081102645809
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 88.14978533333333
$ws.Cells.Item(7, 14).Value = 264.449356
$ws.Cells.Item(7, 15).Value = 0.7247753838328104
$ws.Cells.Item(7, 16).Value = 0.7247753838328105
$ws.Cells.Item(7, 17).Value = 163.4741001165462
$ws.Cells.Item(7, 18).Value = 1471.266901048916
$ws.Cells.Item(7, 19).Value = 0.6568698084878006
$ws.Cells.Item(7, 20).Value = 0.6568698084878007

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Btc"
$ws.Cells.Item(8, 3).Value = "Egfr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.854503666666667
$ws.Cells.Item(8, 8).Value = 5.563511
$ws.Cells.Item(8, 9).Value = 0.9063081102645809
$ws.Cells.Item(8, 10).Value = 0.9063081102645809
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.24063
$ws.Cells.Item(8, 14).Value = 0.72189
$ws.Cells.Item(8, 15).Value = 0.001978481285600361
$ws.Cells.Item(8, 16).Value = 0.001978481285600361
$ws.Cells.Item(8, 17).Value = 0.4462492173100001
$ws.Cells.Item(8, 18).Value = 4.01624295579
$ws.Cells.Item(8, 19).Value = 0.001793113635146302
$ws.Cells.Item(8, 20).Value = 0.001793113635146302

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Btc"
$ws.Cells.Item(9, 3).Value = "Egfr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.854503666666667
$ws.Cells.Item(9, 8).Value = 5.563511
$ws.Cells.Item(9, 9).Value = 0.9063081102645809
$ws.Cells.Item(9, 10).Value = 0.9063081102645809
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 30.53182233333333
$ws.Cells.Item(9, 14).Value = 91.595467
$ws.Cells.Item(9, 15).Value = 0.2510353617660938
$ws.Cells.Item(9, 16).Value = 0.2510353617660938
$ws.Cells.Item(9, 17).Value = 56.6213764671819
$ws.Cells.Item(9, 18).Value = 509.592388204637
$ws.Cells.Item(9, 19).Value = 0.2275153843318139
$ws.Cells.Item(9, 20).Value = 0.2275153843318139
